$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the title heading ---
$titlePara = $d.Paragraphs.First

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r/>' +
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Explore the ancient ruins of El Dorado in Book of Souls II. Enjoy unique features, such as Snake Wilds and two types of free spins for high payout potential. Play now for free.</w:t></w:r>' +
           '</w:p>' +
           '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>'

$insertionPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertionPoint.InsertXML($metaXml)

# InsertXML needed two <w:p> fragments to force the paragraph split at the right spot;
# the 2nd one is an empty spacer paragraph that is only there to carry the split - remove it.
$spacerPara = $titlePara.Next().Next()
$spacerPara.Range.Delete()

# --- 2. Remove the duplicated bold "Play Book of Souls II..." paragraph near the end ---
# (search only the part of the document AFTER the paragraphs we just touched, so we don't
#  re-match the heading/meta-description text at the top of the document)
$old = "Play Book of Souls II: El Dorado for Free - Review"
$searchStart = $d.Paragraphs(3).Range.Start
$dupeRange = $d.Range($searchStart, $d.Content.End)
$dupeRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$dupePara = $d.Range($dupeRange.Start, $dupeRange.End)
$dupePara.Expand(4)
$dupePara.Delete()

# --- 3. Replace the text of the italic "Explore the ancient ruins..." paragraph ---
# (again scope the search to after the top-of-document content so only the bottom, italic
#  occurrence is matched and replaced)
$oldDesc = "Explore the ancient ruins of El Dorado in Book of Souls II. Enjoy unique features, such as Snake Wilds and two types of free spins for high payout potential. Play now for free."
$newDesc = 'Create a feature image fitting the game "Book of Souls II: El Dorado". The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding the Sacred Book of Souls and standing in front of the entrance to the hidden temple. In the background, the temple should be visible with a sense of mystery and adventure.'

$searchStart2 = $d.Paragraphs(3).Range.Start
$descRange = $d.Range($searchStart2, $d.Content.End)
$descRange.Find.Execute($oldDesc, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
# Set .Text directly (instead of passing ReplaceWith to Find.Execute) so straight quotes in the
# new text are not auto-converted into smart quotes.
$descRange.Text = $newDesc
